$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Non-Ferromagnetic" note in column G for the affected BOM rows (new shared string)
$rows = @(8, 9, 10, 11, 12, 13, 17, 18)
foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Non-Ferromagnetic"
}

# Give column G an explicit width (stored width ends up at 19 characters)
$ws.Columns.Item(7).ColumnWidth = 18.14

# Update the view: move the active selection to a single cell, G20
$ws.Range("G20").Select()
